$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update changed values (rows 2-6) ---
$ws.Range("D2").Value = 3399
$ws.Range("E2").Value = -90
$ws.Range("F2").Value = -90
$ws.Range("G2").Value = -193
$ws.Range("H2").Value = -185
$ws.Range("I2").Value = -181
$ws.Range("J2").Value = -4
$ws.Range("K2").Value = 2261
$ws.Range("L2").Value = 1436
$ws.Range("M2").Value = 825
$ws.Range("N2").Value = 720
$ws.Range("O2").Value = 105
$ws.Range("P2").Value = 360
$ws.Range("Q2").Value = -156
$ws.Range("R2").Value = 190
$ws.Range("S2").Value = -125
$ws.Range("T2").Value = 138
$ws.Range("U2").Value = -295
$ws.Range("V2").Value = 1127
$ws.Range("W2").Value = -2.65
$ws.Range("X2").Value = -5.44
$ws.Range("Y2").Value = -24.37
$ws.Range("Z2").Value = -7.15
$ws.Range("AA2").Value = 174.07
$ws.Range("AB2").Value = 106.52
$ws.Range("AC2").Value = -228
$ws.Range("AD2").Value = -2.38
$ws.Range("AE2").Value = 926
$ws.Range("AF2").Value = 0.59
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 79411792
$ws.Range("D3").Value = 2745
$ws.Range("E3").Value = -97
$ws.Range("F3").Value = -97
$ws.Range("G3").Value = -149
$ws.Range("H3").Value = -139
$ws.Range("I3").Value = -133
$ws.Range("J3").Value = -6
$ws.Range("K3").Value = 2228
$ws.Range("L3").Value = 1425
$ws.Range("M3").Value = 803
$ws.Range("N3").Value = 710
$ws.Range("O3").Value = 93
$ws.Range("P3").Value = 417
$ws.Range("Q3").Value = 170
$ws.Range("R3").Value = -119
$ws.Range("S3").Value = 99
$ws.Range("T3").Value = 147
$ws.Range("U3").Value = 23
$ws.Range("V3").Value = 1110
$ws.Range("W3").Value = -3.54
$ws.Range("X3").Value = -5.07
$ws.Range("Y3").Value = -18.56
$ws.Range("Z3").Value = -6.2
$ws.Range("AA3").Value = 177.39
$ws.Range("AB3").Value = 76.48
$ws.Range("AC3").Value = -167
$ws.Range("AD3").Value = -6.32
$ws.Range("AE3").Value = 787
$ws.Range("AF3").Value = 1.34
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 91956911
$ws.Range("D4").Value = 2298
$ws.Range("E4").Value = -35
$ws.Range("F4").Value = -35
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 11
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 2120
$ws.Range("L4").Value = 1147
$ws.Range("M4").Value = 973
$ws.Range("N4").Value = 877
$ws.Range("O4").Value = 96
$ws.Range("P4").Value = 584
$ws.Range("Q4").Value = 47
$ws.Range("R4").Value = 8
$ws.Range("S4").Value = -200
$ws.Range("T4").Value = 20
$ws.Range("U4").Value = 27
$ws.Range("V4").Value = 891
$ws.Range("W4").Value = -1.52
$ws.Range("X4").Value = 0.5
$ws.Range("Y4").Value = 1.03
$ws.Range("Z4").Value = 0.53
$ws.Range("AA4").Value = 117.92
$ws.Range("AB4").Value = 54.58
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 105.68
$ws.Range("AE4").Value = 762
$ws.Range("AF4").Value = 1.17
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 116725477
$ws.Range("D5").Value = 2546
$ws.Range("E5").Value = 69
$ws.Range("F5").Value = 69
$ws.Range("G5").Value = -89
$ws.Range("H5").Value = -103
$ws.Range("I5").Value = -109
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 2517
$ws.Range("L5").Value = 1550
$ws.Range("M5").Value = 967
$ws.Range("N5").Value = 866
$ws.Range("O5").Value = 101
$ws.Range("P5").Value = 584
$ws.Range("Q5").Value = 123
$ws.Range("R5").Value = -156
$ws.Range("S5").Value = 26
$ws.Range("T5").Value = 35
$ws.Range("U5").Value = 87
$ws.Range("V5").Value = 1238
$ws.Range("W5").Value = 2.73
$ws.Range("X5").Value = -4.06
$ws.Range("Y5").Value = -12.5
$ws.Range("Z5").Value = -4.45
$ws.Range("AA5").Value = 160.21
$ws.Range("AB5").Value = 43.47
$ws.Range("AC5").Value = -93
$ws.Range("AD5").Value = -6.94
$ws.Range("AE5").Value = 753
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 116725477
$ws.Range("D6").Value = 2691
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5
$ws.Range("G6").Value = -119
$ws.Range("H6").Value = -112
$ws.Range("I6").Value = -110
$ws.Range("K6").Value = 1933
$ws.Range("L6").Value = 899
$ws.Range("M6").Value = 1034
$ws.Range("N6").Value = 1029
$ws.Range("P6").Value = 792
$ws.Range("Q6").Value = 38
$ws.Range("R6").Value = 86
$ws.Range("S6").Value = -126
$ws.Range("T6").Value = 152
$ws.Range("U6").Value = -114
$ws.Range("V6").Value = 624
$ws.Range("W6").Value = 0.17
$ws.Range("X6").Value = -4.17
$ws.Range("Y6").Value = -11.63
$ws.Range("Z6").Value = -5.05
$ws.Range("AA6").Value = 86.92
$ws.Range("AB6").Value = 20.55
$ws.Range("AC6").Value = -87
$ws.Range("AD6").Value = -8.85
$ws.Range("AE6").Value = 656
$ws.Range("AF6").Value = 1.17
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 158472594

# --- Clear removed cells ---
# Rows 7-9: remove all data columns (D:AJ), keeping only A-C (estimate years cleared)
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

# Row 6: AG6/AH6 (cash DPS / yield) no longer reported - clear them
$ws.Range("AG6:AH6").ClearContents()

Write-Output "Applied IFRS list corrections"